# Natmi following Dr Hou advice
# Rewrite Mmp12-Plaur LR-pair rows: recompute stats for existing
# "M1"/"M2" sending-cluster rows and add a new "Neutro" sending-cluster
# block (rows 14-19), mirroring the same six target clusters
# (ECs, FAPs, M1, M2, Neutro, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Mmp12"
$ws.Range("C2").Value = "Plaur"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.573241999999999
$ws.Range("H2").Value = 28.719726
$ws.Range("I2").Value = 0.443082330403876
$ws.Range("J2").Value = 0.443082330403876
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.4437745
$ws.Range("N2").Value = 28.887549
$ws.Range("O2").Value = 0.1011784119468053
$ws.Range("P2").Value = 0.07108478932534294
$ws.Range("Q2").Value = 138.273748681929
$ws.Range("R2").Value = 829.6424920915739
$ws.Range("S2").Value = 0.04483036655195387
$ws.Range("T2").Value = 0.03149641411054152

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Mmp12"
$ws.Range("C3").Value = "Plaur"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.573241999999999
$ws.Range("H3").Value = 28.719726
$ws.Range("I3").Value = 0.443082330403876
$ws.Range("J3").Value = 0.443082330403876
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.722763333333333
$ws.Range("N3").Value = 20.16829
$ws.Range("O3").Value = 0.04709285083070759
$ws.Range("P3").Value = 0.04962894725691061
$ws.Range("Q3").Value = 64.35864029872666
$ws.Range("R3").Value = 579.2277626885399
$ws.Range("S3").Value = 0.02086601009143203
$ws.Range("T3").Value = 0.021989709606083

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Mmp12"
$ws.Range("C4").Value = "Plaur"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.573241999999999
$ws.Range("H4").Value = 28.719726
$ws.Range("I4").Value = 0.443082330403876
$ws.Range("J4").Value = 0.443082330403876
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 46.33817233333334
$ws.Range("N4").Value = 139.014517
$ws.Range("O4").Value = 0.3245981643651428
$ws.Range("P4").Value = 0.342078784673262
$ws.Range("Q4").Value = 443.6065375847047
$ws.Range("R4").Value = 3992.458838262342
$ws.Range("S4").Value = 0.1438237111117278
$ws.Range("T4").Value = 0.1515690650947547

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Mmp12"
$ws.Range("C5").Value = "Plaur"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.573241999999999
$ws.Range("H5").Value = 28.719726
$ws.Range("I5").Value = 0.443082330403876
$ws.Range("J5").Value = 0.443082330403876
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 54.450333
$ws.Range("N5").Value = 163.350999
$ws.Range("O5").Value = 0.3814237215427815
$ws.Range("P5").Value = 0.4019645747723113
$ws.Range("Q5").Value = 521.2662147895859
$ws.Range("R5").Value = 4691.395933106273
$ws.Range("S5").Value = 0.1690021114124947
$ws.Range("T5").Value = 0.1781034005299187

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Mmp12"
$ws.Range("C6").Value = "Plaur"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.573241999999999
$ws.Range("H6").Value = 28.719726
$ws.Range("I6").Value = 0.443082330403876
$ws.Range("J6").Value = 0.443082330403876
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.359317
$ws.Range("N6").Value = 40.077951
$ws.Range("O6").Value = 0.09358180430980555
$ws.Range("P6").Value = 0.09862147541234521
$ws.Range("Q6").Value = 127.891974595714
$ws.Range("R6").Value = 1151.027771361426
$ws.Range("S6").Value = 0.04146444393698812
$ws.Range("T6").Value = 0.04369743315357048

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Mmp12"
$ws.Range("C7").Value = "Plaur"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.573241999999999
$ws.Range("H7").Value = 28.719726
$ws.Range("I7").Value = 0.443082330403876
$ws.Range("J7").Value = 0.443082330403876
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.441137
$ws.Range("N7").Value = 14.882274
$ws.Range("O7").Value = 0.05212504700475732
$ws.Range("P7").Value = 0.03662142855982794
$ws.Range("Q7").Value = 71.235805256154
$ws.Range("R7").Value = 427.414831536924
$ws.Range("S7").Value = 0.02309568729927945
$ws.Range("T7").Value = 0.01622630790900762

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Mmp12"
$ws.Range("C8").Value = "Plaur"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.980403
$ws.Range("H8").Value = 35.941209
$ws.Range("I8").Value = 0.5544939614414415
$ws.Range("J8").Value = 0.5544939614414415
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.4437745
$ws.Range("N8").Value = 28.887549
$ws.Range("O8").Value = 0.1011784119468053
$ws.Range("P8").Value = 0.07108478932534294
$ws.Range("Q8").Value = 173.0422393511235
$ws.Range("R8").Value = 1038.253436106741
$ws.Range("S8").Value = 0.05610281845273815
$ws.Range("T8").Value = 0.0394160864312397

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Mmp12"
$ws.Range("C9").Value = "Plaur"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.980403
$ws.Range("H9").Value = 35.941209
$ws.Range("I9").Value = 0.5544939614414415
$ws.Range("J9").Value = 0.5544939614414415
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.722763333333333
$ws.Range("N9").Value = 20.16829
$ws.Range("O9").Value = 0.04709285083070759
$ws.Range("P9").Value = 0.04962894725691061
$ws.Range("Q9").Value = 80.54141400695667
$ws.Range("R9").Value = 724.87272606261
$ws.Range("S9").Value = 0.02611270141268993
$ws.Range("T9").Value = 0.02751895156665272

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Mmp12"
$ws.Range("C10").Value = "Plaur"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.980403
$ws.Range("H10").Value = 35.941209
$ws.Range("I10").Value = 0.5544939614414415
$ws.Range("J10").Value = 0.5544939614414415
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 46.33817233333334
$ws.Range("N10").Value = 139.014517
$ws.Range("O10").Value = 0.3245981643651428
$ws.Range("P10").Value = 0.342078784673262
$ws.Range("Q10").Value = 555.1499788367838
$ws.Range("R10").Value = 4996.349809531053
$ws.Range("S10").Value = 0.1799877220354482
$ws.Range("T10").Value = 0.1896806204385509

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Mmp12"
$ws.Range("C11").Value = "Plaur"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 11.980403
$ws.Range("H11").Value = 35.941209
$ws.Range("I11").Value = 0.5544939614414415
$ws.Range("J11").Value = 0.5544939614414415
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 54.450333
$ws.Range("N11").Value = 163.350999
$ws.Range("O11").Value = 0.3814237215427815
$ws.Range("P11").Value = 0.4019645747723113
$ws.Range("Q11").Value = 652.336932824199
$ws.Range("R11").Value = 5871.032395417791
$ws.Range("S11").Value = 0.2114971503459942
$ws.Range("T11").Value = 0.2228869294246234

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Mmp12"
$ws.Range("C12").Value = "Plaur"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 11.980403
$ws.Range("H12").Value = 35.941209
$ws.Range("I12").Value = 0.5544939614414415
$ws.Range("J12").Value = 0.5544939614414415
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 13.359317
$ws.Range("N12").Value = 40.077951
$ws.Range("O12").Value = 0.09358180430980555
$ws.Range("P12").Value = 0.09862147541234521
$ws.Range("Q12").Value = 160.050001464751
$ws.Range("R12").Value = 1440.450013182759
$ws.Range("S12").Value = 0.05189054539058184
$ws.Range("T12").Value = 0.05468501258459101

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Mmp12"
$ws.Range("C13").Value = "Plaur"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 11.980403
$ws.Range("H13").Value = 35.941209
$ws.Range("I13").Value = 0.5544939614414415
$ws.Range("J13").Value = 0.5544939614414415
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 7.441137
$ws.Range("N13").Value = 14.882274
$ws.Range("O13").Value = 0.05212504700475732
$ws.Range("P13").Value = 0.03662142855982794
$ws.Range("Q13").Value = 89.14782003821101
$ws.Range("R13").Value = 534.886920229266
$ws.Range("S13").Value = 0.02890302380398923
$ws.Range("T13").Value = 0.02030636099578374

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Mmp12"
$ws.Range("C14").Value = "Plaur"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05236666666666667
$ws.Range("H14").Value = 0.1571
$ws.Range("I14").Value = 0.002423708154682566
$ws.Range("J14").Value = 0.002423708154682567
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 14.4437745
$ws.Range("N14").Value = 28.887549
$ws.Range("O14").Value = 0.1011784119468053
$ws.Range("P14").Value = 0.07108478932534294
$ws.Range("Q14").Value = 0.75637232465
$ws.Range("R14").Value = 4.5382339479
$ws.Range("S14").Value = 0.0002452269421133041
$ws.Range("T14").Value = 0.000172288783561726

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Mmp12"
$ws.Range("C15").Value = "Plaur"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05236666666666667
$ws.Range("H15").Value = 0.1571
$ws.Range("I15").Value = 0.002423708154682566
$ws.Range("J15").Value = 0.002423708154682567
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 6.722763333333333
$ws.Range("N15").Value = 20.16829
$ws.Range("O15").Value = 0.04709285083070759
$ws.Range("P15").Value = 0.04962894725691061
$ws.Range("Q15").Value = 0.3520487065555556
$ws.Range("R15").Value = 3.168438359
$ws.Range("S15").Value = 0.0001141393265856357
$ws.Range("T15").Value = 0.0001202860841748852

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Mmp12"
$ws.Range("C16").Value = "Plaur"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05236666666666667
$ws.Range("H16").Value = 0.1571
$ws.Range("I16").Value = 0.002423708154682566
$ws.Range("J16").Value = 0.002423708154682567
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 46.33817233333334
$ws.Range("N16").Value = 139.014517
$ws.Range("O16").Value = 0.3245981643651428
$ws.Range("P16").Value = 0.342078784673262
$ws.Range("Q16").Value = 2.426575624522223
$ws.Range("R16").Value = 21.83918062070001
$ws.Range("S16").Value = 0.0007867312179667885
$ws.Range("T16").Value = 0.000829099139956487

# Row 17
$ws.Range("A17").Value = "Neutro"
$ws.Range("B17").Value = "Mmp12"
$ws.Range("C17").Value = "Plaur"
$ws.Range("D17").Value = "M2"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05236666666666667
$ws.Range("H17").Value = 0.1571
$ws.Range("I17").Value = 0.002423708154682566
$ws.Range("J17").Value = 0.002423708154682567
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 54.450333
$ws.Range("N17").Value = 163.350999
$ws.Range("O17").Value = 0.3814237215427815
$ws.Range("P17").Value = 0.4019645747723113
$ws.Range("Q17").Value = 2.8513824381
$ws.Range("R17").Value = 25.6624419429
$ws.Range("S17").Value = 0.0009244597842926121
$ws.Range("T17").Value = 0.0009742448177691613

# Row 18
$ws.Range("A18").Value = "Neutro"
$ws.Range("B18").Value = "Mmp12"
$ws.Range("C18").Value = "Plaur"
$ws.Range("D18").Value = "Neutro"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.05236666666666667
$ws.Range("H18").Value = 0.1571
$ws.Range("I18").Value = 0.002423708154682566
$ws.Range("J18").Value = 0.002423708154682567
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 13.359317
$ws.Range("N18").Value = 40.077951
$ws.Range("O18").Value = 0.09358180430980555
$ws.Range("P18").Value = 0.09862147541234521
$ws.Range("Q18").Value = 0.6995829002333334
$ws.Range("R18").Value = 6.2962461021
$ws.Range("S18").Value = 0.0002268149822355838
$ws.Range("T18").Value = 0.0002390296741837274

# Row 19
$ws.Range("A19").Value = "Neutro"
$ws.Range("B19").Value = "Mmp12"
$ws.Range("C19").Value = "Plaur"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.05236666666666667
$ws.Range("H19").Value = 0.1571
$ws.Range("I19").Value = 0.002423708154682566
$ws.Range("J19").Value = 0.002423708154682567
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 7.441137
$ws.Range("N19").Value = 14.882274
$ws.Range("O19").Value = 0.05212504700475732
$ws.Range("P19").Value = 0.03662142855982794
$ws.Range("Q19").Value = 0.3896675409000001
$ws.Range("R19").Value = 2.3380052454
$ws.Range("S19").Value = 0.0001263359014886424
$ws.Range("T19").Value = 0.00008875965503658003
